$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.717.49"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "2.294.13"
$ws.Range("E3").Value = "  -0.59%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'97.13"
$ws.Range("E5").Value = "  +2.66%  "
$ws.Range("D6").Value = "'269.18"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").Value = "'0.622"
$ws.Range("E7").Value = "  +0.51%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -1.74%  "
$ws.Range("D10").Value = "'45.56"
$ws.Range("E10").Value = "  +2.69%  "
$ws.Range("D11").Value = "'0.0931"
$ws.Range("E11").Value = "  -0.42%  "
$ws.Range("E12").Value = "  -2.10%  "
$ws.Range("E13").Value = "  +1.99%  "
$ws.Range("D14").Value = "2.639.87"
$ws.Range("E14").Value = "  -0.53%  "
$ws.Range("D15").Value = "'15.49"
$ws.Range("E15").Value = "  +0.81%  "
$ws.Range("D16").Value = "'0.851"
$ws.Range("E16").Value = "  -1.18%  "
$ws.Range("D17").Value = "2.290.17"
$ws.Range("E17").Value = "  -0.98%  "
$ws.Range("D18").Value = "43.660.52"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("E19").Value = "  +4.26%  "
$ws.Range("E20").Value = "  -2.24%  "
$ws.Range("D21").Value = "'71.94"
$ws.Range("E21").Value = "  +0.87%  "
$ws.Range("E22").Value = "  +12.19%  "
$ws.Range("D23").Value = "'232.91"
$ws.Range("E23").Value = "  -1.27%  "
$ws.Range("D24").Value = "'9.07"
$ws.Range("E24").Value = "  -4.67%  "
$ws.Range("D25").Value = "'2.72"
$ws.Range("E25").Value = "  +8.89%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").Value = "'11.25"
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("E28").Value = "  +2.31%  "
$ws.Range("D29").Value = "'39.59"
$ws.Range("E29").Value = "  +2.66%  "
$ws.Range("D30").Value = "'2.29"
$ws.Range("E30").Value = "  -1.89%  "
$ws.Range("E31").Value = "  +1.84%  "
$ws.Range("D32").Value = "'21.85"
$ws.Range("E32").Value = "  -2.65%  "
$ws.Range("D33").Value = "'0.0900"
$ws.Range("E33").Value = "  +0.19%  "
$ws.Range("D34").Value = "'5.39"
$ws.Range("E34").Value = "  -2.58%  "
$ws.Range("D35").Value = "'0.125"
$ws.Range("E35").Value = "  -0.51%  "
$ws.Range("D36").Value = "'0.108"
$ws.Range("E36").Value = "  -0.77%  "
$ws.Range("E37").Value = "  -1.64%  "
$ws.Range("D38").Value = "'4.42"
$ws.Range("E38").Value = "  -1.47%  "
$ws.Range("E39").Value = "  -1.85%  "
$ws.Range("D40").Value = "'0.239"
$ws.Range("E40").Value = "  +2.57%  "
$ws.Range("E41").Value = "  +0.56%  "
$ws.Range("D42").Value = "'12.32"
$ws.Range("E42").Value = "  +1.66%  "
$ws.Range("E43").Value = "  -0.72%  "
$ws.Range("D44").Value = "'64.70"
$ws.Range("E44").Value = "  +4.89%  "
$ws.Range("E45").Value = "  -2.74%  "
$ws.Range("D46").Value = "'5.15"
$ws.Range("E46").Value = "  -5.46%  "
$ws.Range("E47").Value = "  -0.41%  "
$ws.Range("D48").Value = "'97.28"
$ws.Range("E48").Value = "  -2.88%  "
$ws.Range("E49").Value = "  -0.85%  "
$ws.Range("D50").Value = "'1.50"
$ws.Range("E50").Value = "  +10.50%  "
$ws.Range("D51").Value = "2.519.98"
$ws.Range("E51").Value = "  -0.47%  "
